$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 0.2476608498783861
    2  = -0.0099999996011277403
    3  = -0.0089999995973677471
    4  = 0.061995794208851152
    5  = -0.0059999996049162618
    6  = -0.0059999995922694893
    7  = -0.019999999516482347
    8  = -0.019999999512418931
    9  = -0.0059999995829072006
    10 = -0.0059999995795863015
    11 = 0.010427393336261304
    12 = -0.0059999995783726057
    13 = -0.0059999995739223877
    14 = -0.011999999540806883
    15 = -0.0059999995720430022
    16 = -0.005999999571629111
    17 = 0.0081517870173097862
    18 = -0.0089999995551810485
    19 = -0.0089999996061917997
    20 = -0.0089999996024303641
    21 = -0.0089999996018610418
    22 = -0.008999999601473796
    23 = -0.0089999995937510846
    24 = -0.041999999411447497
    25 = -0.041999999408152355
    26 = -0.054182501198184241
    27 = -0.0059999995915038795
    28 = -0.0059999995909993942
    29 = -0.011999999559217045
    30 = 0.055191259243931867
    31 = -0.01499999953994724
    32 = -0.020999999507391287
    33 = -0.0059999995870079204
}

foreach ($row in $values.Keys | Sort-Object) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
